# Updates the cryptos list (prices, volume %, and a few reordered
# coin rows) to match the latest scrape.
#
# For cells in column D whose new value looks like a plain decimal
# number (e.g. "608.65"), we briefly force the cell's NumberFormat to
# Text ("@") before assigning the value so Excel stores the exact
# string instead of silently converting it to a floating point number
# (which would introduce binary rounding noise such as
# 608.64999999999998). Afterwards we reset the cell style back to
# "Normal" so no stray number-format/style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.903.26'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '2.698.72'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.124'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000201'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.18%  '
$ws.Range("D15").Value = '3.191.50'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '65.728.74'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '2.690.22'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '359.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.36%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000106'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.169'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("E30").Value = '  +3.58%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '534.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.432'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.11%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0612'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.91%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0266'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.07%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.658'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0985'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.68%  '
